# Remove the "Lương" sheet entirely (salary report no longer generated here).
$wb = $excel.ActiveWorkbook
$salarySheet = $wb.Worksheets.Item("Lương")
$salarySheet.Delete()

$ws = $wb.Worksheets.Item("Đơn sale chính")

# Insert a new column before the current column G ("Tên dịch vụ"),
# shifting G:V to H:W, to make room for "Nhóm dịch vụ".
$ws.Columns.Item(7).Insert()

# Insert a new row before the current row 3 ("Tổng"), shifting it to row 4,
# so we can add a second data row.
$ws.Rows.Item(3).Insert()

# --- Header row (row 1) ---
$ws.Cells.Item(1,7).Value  = "Nhóm dịch vụ"
$ws.Cells.Item(1,24).Value = "Tỉ lệ chiết khấu sale chính"
$ws.Cells.Item(1,25).Value = "Tỉ lệ chiết khấu sale phụ"
$ws.Cells.Item(1,26).Value = "Chiết khấu sale chính"
$ws.Cells.Item(1,27).Value = "Chiết khấu sale phụ"

# --- Row 2 (existing order, HD-LUXURY-526) ---
$ws.Cells.Item(2,7).Value  = "Vùng mắt"
$ws.Cells.Item(2,24).Value = 0
$ws.Cells.Item(2,25).Value = 0.02
$ws.Cells.Item(2,26).Value = 0
$ws.Cells.Item(2,27).Value = 20000

# --- Row 3 (new order, HD-LUXURY-535) ---
$ws.Cells.Item(3,1).Value  = "HD-LUXURY"
$ws.Cells.Item(3,2).Value  = 535
# Force text so the dd-mm-yyyy-looking string isn't auto-converted to a date serial.
$ws.Cells.Item(3,3).NumberFormat = "@"
$ws.Cells.Item(3,3).Value  = "07-09-2024"
$ws.Cells.Item(3,4).Value  = "CẦN THƠ"
$ws.Cells.Item(3,5).Value  = "Phan Minh Nguyệt"
$ws.Cells.Item(3,6).Value  = "CTV"
$ws.Cells.Item(3,7).Value  = "Vùng mắt"
$ws.Cells.Item(3,8).Value  = "Phun mày"
$ws.Cells.Item(3,9).Value  = "CTV Ngoài"
$ws.Cells.Item(3,10).Value = 500000
$ws.Cells.Item(3,11).Value = "Đỗ Thị Huyền Trân"
$ws.Cells.Item(3,12).Value = 400000
$ws.Cells.Item(3,13).Value = 900000
$ws.Cells.Item(3,14).Value = 900000
$ws.Cells.Item(3,15).Value = 0
$ws.Cells.Item(3,16).Value = 900000
$ws.Cells.Item(3,17).Value = 0
$ws.Cells.Item(3,18).Value = "Nguyễn Hoàng Yến Quyên"
# Bác sĩ 2 / Phụ phẫu 1-2 / Công phụ phẫu 1-2 are blank numeric cells for this row.
$ws.Cells.Item(3,19).Value = 0
$ws.Cells.Item(3,20).Value = 0
$ws.Cells.Item(3,21).Value = 0
$ws.Cells.Item(3,22).Value = 0
$ws.Cells.Item(3,23).Value = 0
$ws.Cells.Item(3,24).Value = 0
$ws.Cells.Item(3,25).Value = 0.02
$ws.Cells.Item(3,26).Value = 0
$ws.Cells.Item(3,27).Value = 18000

# --- Row 4 ("Tổng" totals row, pushed down from row 3) ---
$ws.Cells.Item(4,2).Value  = 2
$ws.Cells.Item(4,10).Value = 1000000
$ws.Cells.Item(4,12).Value = 1400000
$ws.Cells.Item(4,13).Value = 2400000
$ws.Cells.Item(4,14).Value = 1900000
$ws.Cells.Item(4,15).Value = 0
$ws.Cells.Item(4,16).Value = 1900000
$ws.Cells.Item(4,17).Value = 500000
$ws.Cells.Item(4,24).Value = 0
$ws.Cells.Item(4,25).Value = 0.04
$ws.Cells.Item(4,26).Value = 0
$ws.Cells.Item(4,27).Value = 38000
